# "Data source corrected and updated"
# J1/K1 previously held shared-string placeholders ("r"/"s"); the real
# numeric data source now supplies 0.6 for both. J2:J51 were an earlier
# (wrong) constant of 0.5 and are corrected to match K2:K51's 0.6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace the text placeholders in J1/K1 with the real numeric values.
$ws.Range("J1").Value = 0.6
$ws.Range("K1").Value = 0.6

# Rows 2-51: correct column J from 0.5 to 0.6 (column K already correct).
$ws.Range("J2:J51").Value = 0.6

# Refresh the view: zoom to 100% and move the selection/active cell to K1.
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 39
$ws.Range("K1:K51").Select()
